# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# worksheet with refreshed values, matching the upstream GitHub Actions
# data-refresh commit. Numeric-looking price strings are entered with a
# leading apostrophe so Excel stores them as text (matching the existing
# text-formatted Price column) instead of auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.356.86"
$ws.Range("E2").Value = "  +0.60%  "

$ws.Range("D3").Value = "2.011.33"
$ws.Range("E3").Value = "  +0.29%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'258.50"
$ws.Range("E5").Value = "  +4.75%  "

$ws.Range("D6").Value = "'0.611"
$ws.Range("E6").Value = "  -1.84%  "

$ws.Range("D8").Value = "'56.09"
$ws.Range("E8").Value = "  -6.30%  "

$ws.Range("D9").Value = "'0.386"
$ws.Range("E9").Value = "  +0.09%  "

$ws.Range("E10").Value = "  -5.00%  "

$ws.Range("E11").Value = "  -1.67%  "

$ws.Range("D12").Value = "2.307.55"
$ws.Range("E12").Value = "  +0.28%  "

$ws.Range("D13").Value = "'14.23"
$ws.Range("E13").Value = "  -5.50%  "

$ws.Range("E14").Value = "  -5.20%  "

$ws.Range("D15").Value = "'20.95"
$ws.Range("E15").Value = "  -6.13%  "

$ws.Range("E16").Value = "  -3.94%  "

$ws.Range("D17").Value = "2.031.70"
$ws.Range("E17").Value = "  +1.91%  "

$ws.Range("D18").Value = "37.207.23"
$ws.Range("E18").Value = "  +0.07%  "

$ws.Range("D19").Value = "'69.63"
$ws.Range("E19").Value = "  -0.86%  "

$ws.Range("E20").Value = "  -3.65%  "

$ws.Range("E21").Value = "  -0.88%  "

$ws.Range("D22").Value = "'228.48"
$ws.Range("E22").Value = "  -0.69%  "

$ws.Range("D23").Value = "'2.63"
$ws.Range("E23").Value = "  +6.05%  "

$ws.Range("E24").Value = "  +0.07%  "

$ws.Range("E25").Value = "  -0.05%  "

$ws.Range("D26").Value = "'164.70"
$ws.Range("E26").Value = "  +0.08%  "

$ws.Range("D27").Value = "'8.93"
$ws.Range("E27").Value = "  -5.30%  "

$ws.Range("D28").Value = "'19.66"
$ws.Range("E28").Value = "  +0.07%  "

$ws.Range("D29").Value = "'0.129"
$ws.Range("E29").Value = "  -6.63%  "

$ws.Range("E30").Value = "  -4.91%  "

$ws.Range("E31").Value = "  -1.19%  "

$ws.Range("E32").Value = "  -3.40%  "

$ws.Range("D33").Value = "'0.0646"
$ws.Range("E33").Value = "  -1.35%  "

$ws.Range("E34").Value = "  +1.95%  "

$ws.Range("E35").Value = "  -2.07%  "

$ws.Range("E36").Value = "  +0.86%  "

$ws.Range("E37").Value = "  -1.56%  "

$ws.Range("E38").Value = "  -0.07%  "

$ws.Range("D39").Value = "'5.25"
$ws.Range("E39").Value = "  -2.14%  "

$ws.Range("D40").Value = "'3.04"
$ws.Range("E40").Value = "  +3.93%  "

$ws.Range("D41").Value = "'1.20"
$ws.Range("E41").Value = "  +1.61%  "

$ws.Range("D42").Value = "'0.0933"
$ws.Range("E42").Value = "  -5.14%  "

$ws.Range("E43").Value = "  -0.86%  "

$ws.Range("D44").Value = "1.407.55"
$ws.Range("E44").Value = "  +2.69%  "

$ws.Range("D45").Value = "'89.90"
$ws.Range("E45").Value = "  -2.04%  "

$ws.Range("D46").Value = "'15.69"
$ws.Range("E46").Value = "  -5.67%  "

$ws.Range("E47").Value = "  -2.57%  "

$ws.Range("D48").Value = "'7.05"
$ws.Range("E48").Value = "  -4.52%  "

$ws.Range("E49").Value = "  +2.55%  "

$ws.Range("D50").Value = "2.199.86"
$ws.Range("E50").Value = "  +0.31%  "

$ws.Range("E51").Value = "  -6.66%  "
